$d = $word.ActiveDocument

# Helper: convert an RRGGBB hex string into the BGR-packed integer Word's
# Font.Color property expects (same packing as the RGB() VBA macro).
function Get-WordColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$redColor = Get-WordColor "C9211E"

# Find `term` as a whole word inside the given paragraph (by 1-based index)
# and paint it with `color`.
function Set-TermColor($paragraphIndex, $term, $color) {
    $p = $d.Paragraphs($paragraphIndex)
    $searchRange = $p.Range
    $found = $searchRange.Find.Execute($term, $true, $true, $false, $false, $false, `
                                        $true, 1, $false, "", 0)
    if ($found) {
        $searchRange.Font.Color = $color
    }
}

# "Families to do:" section -> "Continue Saxifragaceae, Cucurbitaceae"
# Highlight the two family names that are the immediate next steps.
Set-TermColor 43 "Saxifragaceae" $redColor
Set-TermColor 43 "Cucurbitaceae" $redColor

# "Add these:" list -> "Anisophylleaceae, Apodanthaceae, Corynocarpaceae,
# Coriariaceae, Tetramelaceae, Datiscaceae, Begoniaceae"
# Highlight only the families that should be tackled next.
Set-TermColor 45 "Apodanthaceae" $redColor
Set-TermColor 45 "Datiscaceae" $redColor
Set-TermColor 45 "Begoniaceae" $redColor

# "Rosaceae, Rhamnaceae, Elaeagnaceae, ..." -> highlight Rosaceae.
Set-TermColor 46 "Rosaceae" $redColor
